$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing GPA values (D3:D5) ---
$ws.Range("D3").Value = 3.97
$ws.Range("D4").Value = 3.23
$ws.Range("D5").Value = 3.55

# --- Column A (names) for the new rows, in row order ---
$ws.Range("A6").Value  = "Jason"
$ws.Range("A7").Value  = "Coco"
$ws.Range("A8").Value  = "Rose"
$ws.Range("A9").Value  = "Tony"
$ws.Range("A10").Value = "Avery"
$ws.Range("A11").Value = "Jessica"
$ws.Range("A12").Value = "Jeremy"
$ws.Range("A13").Value = "Eve"

# --- Column B (IDs) ---
$ws.Range("B10").Value = "m123"
$ws.Range("B11").Value = "m456"
$ws.Range("B12").Value = "n789"
$ws.Range("B13").Value = "n135"
$ws.Range("B9").Value  = "l123135"
$ws.Range("B8").Value  = "l712389"
$ws.Range("B7").Value  = "l45426"
$ws.Range("B6").Value  = "l12323"

# --- Column C (emails) + hyperlinks ---
$ws.Range("C6").Value  = "l12323@qq.com"
$ws.Range("C7").Value  = "l45426@qq.com"
$ws.Range("C8").Value  = "l712389@qq.com"
$ws.Range("C9").Value  = "l123135@qq.com"
$ws.Range("C10").Value = "m123@qq.com"
$ws.Range("C11").Value = "m456@qq.com"
$ws.Range("C12").Value = "n789@qq.com"
$ws.Range("C13").Value = "n135@qq.com"

$ws.Hyperlinks.Add($ws.Range("C6"),  "mailto:l12323@qq.com")
$ws.Hyperlinks.Add($ws.Range("C7"),  "mailto:l45426@qq.com")
$ws.Hyperlinks.Add($ws.Range("C8"),  "mailto:l712389@qq.com")
$ws.Hyperlinks.Add($ws.Range("C9"),  "mailto:l123135@qq.com")
$ws.Hyperlinks.Add($ws.Range("C10"), "mailto:m123@qq.com")
$ws.Hyperlinks.Add($ws.Range("C11"), "mailto:m456@qq.com")
$ws.Hyperlinks.Add($ws.Range("C12"), "mailto:n789@qq.com")
$ws.Hyperlinks.Add($ws.Range("C13"), "mailto:n135@qq.com")

# Match the hyperlink style used by the existing hyperlink cells (C2:C5)
$ws.Range("C6:C13").Style = $ws.Range("C2").Style

# --- Column D (GPA-like score) for the new rows ---
$ws.Range("D6").Value  = 2.89
$ws.Range("D7").Value  = 2.64
$ws.Range("D8").Value  = 3.31
$ws.Range("D9").Value  = 3.19
$ws.Range("D10").Value = 3.6
$ws.Range("D11").Value = 3.36
$ws.Range("D12").Value = 3.57
$ws.Range("D13").Value = 3.08

# --- Selection moves to F8 ---
[void]$ws.Range("F8").Select()
